$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Long text blocks (read from here-strings to avoid escaping issues) ---
$objetivosPt = @'
Apresentar os fundamentos do controle automático de processos. Apresentar os diversos instrumentos utilizados para medição de pressão, temperatura, nível e vazão, e os dispositivos dos sistemas de automação de processos.
'@

$resumoPt = @'
Fundamentos do controle automático de processos. Simbologia. Medição de pressão. Medição de temperatura.Medição de nível. Medição de vazão. Controladores. Elementos finais de controle. Transmissão. Sistemas de automação de processos industriais.
'@

$programaPt = @'
1- INTRODUÇÃO AO CONTROLE AUTOMÁTICO DE PROCESSOS. - Introdução. - Malhas aberta e fechada. - Simbologia. - Realimentação. - Diagramas de blocos. - Função de transferência. - Noções de resposta ao degrau do processo de primeira ordem.
2- MEDIÇÃO DE PRESSÃO. - Métodos a coluna líquida. - Métodos a elementos elásticos. - Métodos elétricos.
3- MEDIÇÃO DE TEMPERATURA. - Termopares. - Termômetros de resistência. - Termômetros a volume repleto. - Termômetros bimetálicos. -Termômetros de líquido em vidro. - Pirômetros.
4- MEDIÇÃO DE NÍVEL. - Bóias. - Corpo imerso. - Dispositivos a pressão.
5- MEDIÇÃO DE VAZÃO. - Placa de orifício, Venturi, rotâmetro. - Medidor magnético. - Medidor de vazão mássica Coriolis.
6- CONTROLADORES. - Controle a duas posições, controle proporcional, controle proporcional e integral, e controle proporcional, integral e derivativo. - Controlador Lógico Programável.
7- ELEMENTO FINAL DE CONTROLE. - Válvulas de controle.
8- TRANSMISSÃO. - Transmissão de sinais.
9- SISTEMAS DE AUTOMAÇÃO DE PROCESSOS INDUSTRIAIS.
'@

$biblioPt = @'
1) ALVES, J. L. L. Instrumentação, Controle e Automação de Processos. 2ª ed. Rio de Janeiro: LTC, 2010.
2) BEGA, E. A. (Organizador) Instrumentação Industrial. 3ª ed. Rio de Janeiro: Interciência: IBP, 2011.
3) BALBINOT, A.; BRUSAMARELLO, V. J. Instrumentação e Fundamentos de Medidas. vols 1 e 2. 2ª ed. Rio de Janeiro: LTC, 2011.
4) PERRY, R. H.; CHILTON, C. H. Manual de Engenharia Química. 5ª ed. Rio de Janeiro: Guanabara Dois, 1986. Seção 22.
5) SIGHIERI, L.; NISHINARI, A. Controle Automático de Processos Industriais:      Instrumentação. 2ª ed. São Paulo: Edgard Blücher, 1973.
'@

$shortSyllabusEn = @'
Introduction to automatic process control. Symbology. Pressure measurement. Temperature measurement. Level measurement. Flow measurement. Controllers. Final control element. Transmission. Systems automation of industrial processes.
'@

$longSyllabusEn = @'
1- INTRODUCTION TO AUTOMATIC PROCESS CONTROL. - Introduction. - Open-loop system and closed-loop system. - Symbology. - Feedback control. - Block diagrams. - Transfer function. - Notions of step response of first order process.
2- PRESSURE MEASUREMENT. - Liquid-column methods.  Elastic element methods. - Electrical methods.
3- TEMPERATURE MEASUREMENT. - Thermocouples. - Resistance thermometers. - Filled-system thermometers. - Bimetal thermometers. - Liquid-in-glass thermometers. - Pyrometers.
4- LEVEL MEASUREMENT. - Float-actuated devices. - Pressure devices.
5- FLOW MEASUREMENT. - Orifice meter, Venturi meter, rotameter. - Magnetic flowmeters. - Coriolis mass flowmeters.
6- CONTROLLERS. - On/off control, proportional control, proportional-plus- integral control, proportional-plus-integral-plus-derivative control. - Programmable logic controller.
7- FINAL CONTROL ELEMENT. - Control valves.
8- TRANSMISSION. - Signal transmission.
9- SYSTEMS INDUSTRIAL PROCESS AUTOMATION.
'@

$lob1006 = @'
LOB1006 -  Cálculo IV  (Requisito fraco)

'@

$loq4083 = @'
LOQ4083 -  Fenômenos de Transporte I  (Requisito fraco)

'@

# --- Row 10: Objetivos (PT) body text replaced ---
$ws.Range("B10:C10").Value = $objetivosPt

# --- Row 13: "Programa resumido:" label removed from A13; B13/C13 become the docente name ---
$ws.Range("A13").Clear()
$ws.Range("B13:C13").Value = "5840643 - Luiz Carlos de Queiroz"
$ws.Rows(13).AutoFit()

# --- Row 14: becomes "Programa resumido:" + PT summary text ---
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14:C14").Value = $resumoPt

# --- Row 15: becomes "Short syllabus:" + EN short text ---
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15:C15").Value = $shortSyllabusEn
$ws.Rows(15).RowHeight = 60

# --- Row 16: becomes "Programa:" + PT long syllabus text ---
$ws.Range("A16").Value = "Programa:"
$ws.Range("B16:C16").Value = $programaPt

# --- Row 17: becomes "Syllabus:" + EN long syllabus text (new cells) ---
$ws.Range("A17").Value = "Syllabus:"
$ws.Range("B17:C17").Value = $longSyllabusEn
$ws.Range("B17").Font.Bold = $false
$ws.Range("B17").WrapText = $true
$ws.Range("B17").VerticalAlignment = -4160
$ws.Rows(17).RowHeight = 120

# --- Row 18: becomes "Avaliação:" ; clear the docente name that used to live here ---
$ws.Range("A18").Value = "Avaliação:"
$ws.Range("B18:C18").Clear()
$ws.Rows(18).AutoFit()

# --- Row 19: becomes "Método:" ---
$ws.Range("A19").Value = "Método:"

# --- Row 20: becomes "Critério:" ---
$ws.Range("A20").Value = "Critério:"

# --- Row 21: becomes "Norma de recuperação:" ---
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Rows(21).RowHeight = 60

# --- Row 22: becomes "Bibliografia:" + new bibliography text ---
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22:C22").Value = $biblioPt
$ws.Range("B22").Font.Bold = $false
$ws.Range("B22").WrapText = $true
$ws.Range("B22").VerticalAlignment = -4160
$ws.Rows(22).RowHeight = 120

# --- Row 23: becomes "Requisitos:" only; clear old LOB1006 text that lived here ---
$ws.Range("B23:C23").Clear()
$ws.Range("A23").Value = "Requisitos:"
$ws.Rows(23).AutoFit()

# --- Row 24: now holds the LOB1006 requirement text ---
$ws.Range("B24:C24").Value = $lob1006

# --- Row 25 (new): holds the LOQ4083 requirement text ---
$ws.Range("B25:C25").Value = $loq4083
$ws.Range("B25").Font.Bold = $false
$ws.Range("B25").WrapText = $true
$ws.Range("B25").VerticalAlignment = -4160
$ws.Rows(25).RowHeight = 30

# --- Column layout: split the A:B merged width definition so column A only spans itself ---
$ws.Columns("B").ColumnWidth = $ws.Columns("B").ColumnWidth

Write-Output "edit complete"
